# Sprint 3 planning doc touch-up:
#   - title "Sprint 1 Planning" -> "Sprint 3 Planning"
#   - burndown date "March 3rd, 2022" -> "March 30th, 2022"
#     (the "th" suffix must stay superscripted, same as "rd" was)

$d = $word.ActiveDocument

# 1. Title: bump the sprint number.
$d.Content.Find.Execute("Sprint 1 ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Sprint 3 ", 2)

# 2. Date: change the day-of-month first, while it is still plain bold text
#    (this leaves the separately-formatted superscript run untouched)...
$d.Content.Find.Execute("March 3", $true, $false, $false, $false, $false,
                         $true, 1, $false, "March 30", 2)

# ...then update just the superscripted ordinal suffix "rd" -> "th" so the
# new "30th" keeps the superscript formatting the "rd" run already had.
$d.Content.Find.Execute("rd", $true, $false, $false, $false, $false,
                         $true, 1, $false, "th", 2)
